$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update URL path strings to include path parameters
$ws.Range("A2").Value = "notes/search_title/{title}"
$ws.Range("A3").Value = "notes/search_label/{label}"

# Update the active selection on the sheet
$ws.Range("A5").Select()
